# Updates cryptos price (D) and 1h volume change (E) columns
# per the latest scrape, preserving original text cell typing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.874.91"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.934.67"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.14%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.34"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.34"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.94%  "

$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("E8").Value = "  +4.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.931.93"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.39%  "

$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("E12").Value = "  +3.82%  "

$ws.Range("E13").Value = "  +4.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.86"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.122"
$ws.Range("D15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.419.54"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.87"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +7.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.929.90"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.845.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "416.80"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.92%  "

$ws.Range("E22").Value = "  +7.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.44"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +7.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.98"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.45"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.31%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  +0.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.02"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.41"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.37"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.97"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.10%  "

$ws.Range("E33").Value = "  +4.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.66"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.936"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.17%  "

$ws.Range("E36").Value = "  +6.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0695"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +13.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.28"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.71"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.20%  "

$ws.Range("E40").Value = "  +10.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.107"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "375.59"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.35%  "

$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.694.88"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.13%  "

$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.66"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.236"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.63%  "

$ws.Range("E48").Value = "  +2.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.96"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.95"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.99"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.19%  "

